$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: a new price observation (date 2021-09-10 / serial 44449)
# was inserted ahead of the existing "Terminal La Palmera de La Serena - Ajo"
# series, pushing the rest of the rows (previously 125-139) down by one
# (now 126-140).
$ws.Rows("125:125").Insert()

$ws.Range("A125").Value = 8
$ws.Range("B125").Value = "Terminal La Palmera de La Serena"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44449
$ws.Range("E125").Value = 4
$ws.Range("F125").Value = 100112003
$ws.Range("G125").Value = "Ajo"
$ws.Range("H125").Value = "Chino"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 720
$ws.Range("K125").Value = 15000
$ws.Range("L125").Value = 16000
$ws.Range("M125").Value = 15500
$ws.Range("N125").Value = "`$/caja 10 kilos"
$ws.Range("O125").Value = "China"
$ws.Range("P125").Value = 1550
$ws.Range("Q125").Value = 10
$ws.Range("R125").Value = "Hortaliza"
